$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 6062367.5
$ws.Range("I33").Value = 2840.5
$ws.Range("J33").Value = 12987542
$ws.Range("K33").Value = 2840.5
$ws.Range("L33").Value = 12987542
$ws.Range("M33").Value = -2611.5
$ws.Range("N33").Value = -12988000
$ws.Range("H70").Value = 3059.2942
$ws.Range("I70").Value = 3400.2856
$ws.Range("J70").Value = 2820.6
$ws.Range("K70").Value = 10200.8568
$ws.Range("L70").Value = 8461.799999999999
$ws.Range("M70").Value = -9930.856800000001
$ws.Range("N70").Value = -9001.799999999999
$ws.Range("H73").Value = 3059.2942
$ws.Range("I73").Value = 3400.2856
$ws.Range("J73").Value = 2820.6
$ws.Range("K73").Value = 10200.8568
$ws.Range("L73").Value = 8461.799999999999
$ws.Range("M73").Value = -9264.856800000001
$ws.Range("N73").Value = -10333.8
$ws.Range("H112").Value = 20409120
$ws.Range("J112").Value = 21979020
$ws.Range("L112").Value = 65937060
$ws.Range("N112").Value = -65939276
$ws.Range("H131").Value = 3005
$ws.Range("I131").Value = 489
$ws.Range("J131").Value = 6150
$ws.Range("K131").Value = 1467
$ws.Range("L131").Value = 18450
$ws.Range("M131").Value = 3573
$ws.Range("N131").Value = -28530
$ws.Range("H132").Value = 1859.0605
$ws.Range("I132").Value = 1421.8334
$ws.Range("J132").Value = 6231.3335
$ws.Range("K132").Value = 4265.5002
$ws.Range("L132").Value = 18694.0005
$ws.Range("M132").Value = -1735.5002
$ws.Range("N132").Value = -23754.0005
$ws.Range("H138").Value = 1986.0944
$ws.Range("I138").Value = 1209.6666
$ws.Range("J138").Value = 2998.8262
$ws.Range("K138").Value = 3628.9998
$ws.Range("L138").Value = 8996.4786
$ws.Range("M138").Value = 1511.0002
$ws.Range("N138").Value = -19276.4786
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 611.2708
$ws.Range("I2").Value = 485.41464
$ws.Range("J2").Value = 1348.4286
$ws.Range("K2").Value = 485.41464
$ws.Range("L2").Value = 1348.4286
$ws.Range("M2").Value = -372.41464
$ws.Range("N2").Value = -1574.4286
$ws.Range("H88").Value = 2373
$ws.Range("I88").Value = 2239.2
$ws.Range("J88").Value = 2484.5
$ws.Range("K88").Value = 2239.2
$ws.Range("L88").Value = 2484.5
$ws.Range("M88").Value = -1833.2
$ws.Range("N88").Value = -3296.5
$ws.Range("H91").Value = 2373
$ws.Range("I91").Value = 2239.2
$ws.Range("J91").Value = 2484.5
$ws.Range("K91").Value = 2239.2
$ws.Range("L91").Value = 2484.5
$ws.Range("M91").Value = -835.1999999999998
$ws.Range("N91").Value = -5292.5
$ws.Range("H110").Value = 2028.125
$ws.Range("I110").Value = 1751.8422
$ws.Range("J110").Value = 3078
$ws.Range("K110").Value = 1751.8422
$ws.Range("L110").Value = 3078
$ws.Range("M110").Value = 293.1578
$ws.Range("N110").Value = -7168
$ws.Range("H116").Value = 611.2708
$ws.Range("I116").Value = 485.41464
$ws.Range("J116").Value = 1348.4286
$ws.Range("K116").Value = 485.41464
$ws.Range("L116").Value = 1348.4286
$ws.Range("M116").Value = 1808.58536
$ws.Range("N116").Value = -5936.4286
$ws.Range("H132").Value = 2706814
$ws.Range("I132").Value = 3192.1428
$ws.Range("K132").Value = 9576.428400000001
$ws.Range("M132").Value = -7046.428400000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 611.2708
$ws.Range("I3").Value = 485.41464
$ws.Range("J3").Value = 1348.4286
$ws.Range("K3").Value = 485.41464
$ws.Range("L3").Value = 1348.4286
$ws.Range("M3").Value = -371.41464
$ws.Range("N3").Value = -1576.4286
$ws.Range("H25").Value = 1497.4
$ws.Range("I25").Value = 371.75
$ws.Range("J25").Value = 6000
$ws.Range("K25").Value = 371.75
$ws.Range("L25").Value = 6000
$ws.Range("M25").Value = -136.75
$ws.Range("N25").Value = -6470
$ws.Range("H86").Value = 1955.5294
$ws.Range("I86").Value = 1825.9231
$ws.Range("K86").Value = 1825.9231
$ws.Range("M86").Value = -702.9231
$ws.Range("H89").Value = 1955.5294
$ws.Range("I89").Value = 1825.9231
$ws.Range("K89").Value = 9129.6155
$ws.Range("M89").Value = -3513.6155
$ws.Range("H124").Value = 9912.5
$ws.Range("I124").Value = 9912.5
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 9912.5
$ws.Range("L124").Value = 0
$ws.Range("M124").Value = -5002.5
$ws.Range("N124").ClearContents()
$ws.Range("H134").Value = 35004.12
$ws.Range("I134").Value = 5076.6206
$ws.Range("J134").Value = 251978.5
$ws.Range("K134").Value = 15229.8618
$ws.Range("L134").Value = 755935.5
$ws.Range("M134").Value = -12694.8618
$ws.Range("N134").Value = -761005.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3044.739
$ws.Range("I105").Value = 3008.1365
$ws.Range("K105").Value = 3008.1365
$ws.Range("M105").Value = -1261.1365
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 3714.6667
$ws.Range("J88").Value = 3714.6667
$ws.Range("L88").Value = 11144.0001
$ws.Range("N88").Value = -12000.0001
$ws.Range("H91").Value = 3714.6667
$ws.Range("J91").Value = 3714.6667
$ws.Range("L91").Value = 11144.0001
$ws.Range("N91").Value = -14108.0001
$ws.Range("H122").Value = 2951.2083
$ws.Range("I122").Value = 372.04544
$ws.Range("J122").Value = 5133.577
$ws.Range("K122").Value = 3348.40896
$ws.Range("L122").Value = 46202.193
$ws.Range("M122").Value = -898.4089599999998
$ws.Range("N122").Value = -51102.193
$ws.Range("H137").Value = 10354.35
$ws.Range("I137").Value = 6172.5
$ws.Range("J137").Value = 16627.125
$ws.Range("K137").Value = 18517.5
$ws.Range("L137").Value = 49881.375
$ws.Range("M137").Value = -13417.5
$ws.Range("N137").Value = -60081.375
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 38462924
$ws.Range("I113").Value = 76924050
$ws.Range("K113").Value = 76924050
$ws.Range("M113").Value = -76921880
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 37039308
$ws.Range("I40").Value = 47621340
$ws.Range("J40").Value = 2200.8333
$ws.Range("K40").Value = 47621340
$ws.Range("L40").Value = 2200.8333
$ws.Range("M40").Value = -47621204
$ws.Range("N40").Value = -2472.8333
$ws.Range("H82").Value = 1114610.2
$ws.Range("I82").Value = 3336633.2
$ws.Range("J82").Value = 162314.72
$ws.Range("K82").Value = 3336633.2
$ws.Range("L82").Value = 162314.72
$ws.Range("M82").Value = -3336272.2
$ws.Range("N82").Value = -163036.72
$ws.Range("H85").Value = 1114610.2
$ws.Range("I85").Value = 3336633.2
$ws.Range("J85").Value = 162314.72
$ws.Range("K85").Value = 3336633.2
$ws.Range("L85").Value = 162314.72
$ws.Range("M85").Value = -3335385.2
$ws.Range("N85").Value = -164810.72
$ws.Range("H136").Value = 6220.719
$ws.Range("I136").Value = 5091.343
$ws.Range("J136").Value = 8017.4546
$ws.Range("K136").Value = 15274.029
$ws.Range("L136").Value = 24052.3638
$ws.Range("M136").Value = -12724.029
$ws.Range("N136").Value = -29152.3638
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2600.2
$ws.Range("I81").Value = 2000.3334
$ws.Range("J81").Value = 3500
$ws.Range("K81").Value = 4000.6668
$ws.Range("L81").Value = 7000
$ws.Range("M81").Value = -2939.6668
$ws.Range("N81").Value = -9122
$ws.Range("H84").Value = 2600.2
$ws.Range("I84").Value = 2000.3334
$ws.Range("J84").Value = 3500
$ws.Range("K84").Value = 20003.334
$ws.Range("L84").Value = 35000
$ws.Range("M84").Value = -14699.334
$ws.Range("N84").Value = -45608
$ws.Range("H107").Value = 166667360
$ws.Range("I107").Value = 200000580
$ws.Range("J107").Value = 1280
$ws.Range("K107").Value = 600001740
$ws.Range("L107").Value = 3840
$ws.Range("M107").Value = -599999820
$ws.Range("N107").Value = -7680
$ws.Range("H113").Value = 1100.0667
$ws.Range("I113").Value = 798.30554
$ws.Range("J113").Value = 2307.111
$ws.Range("K113").Value = 2394.91662
$ws.Range("L113").Value = 6921.333
$ws.Range("M113").Value = -224.91662
$ws.Range("N113").Value = -11261.333
$ws.Range("H126").Value = 650.36664
$ws.Range("I126").Value = 480.80768
$ws.Range("J126").Value = 1752.5
$ws.Range("K126").Value = 1442.42304
$ws.Range("L126").Value = 5257.5
$ws.Range("M126").Value = 1027.57696
$ws.Range("N126").Value = -10197.5
$ws.Range("H132").Value = 1655.2703
$ws.Range("I132").Value = 1144.5358
$ws.Range("J132").Value = 3244.2222
$ws.Range("K132").Value = 3433.6074
$ws.Range("L132").Value = 9732.6666
$ws.Range("M132").Value = -903.6074000000003
$ws.Range("N132").Value = -14792.6666
